# Adds a new "test_xlr_n_percent" column (K) to the worksheet/table, containing
# "n (pct%)" style labels for data rows 3-34 (one per table row), matching the
# build_multiple_response / seen_but_answered "n (percent)" summary column.
#
# Percent values are precomputed (i/32*100 rounded half-to-even, i.e. R's / Excel
# display rounding) rather than computed at runtime, to avoid relying on this
# host's Math.Round (which rounds half away from zero, not half-to-even).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "1 (3%)", "2 (6%)", "3 (9%)", "4 (12%)", "5 (16%)", "6 (19%)",
    "7 (22%)", "8 (25%)", "9 (28%)", "10 (31%)", "11 (34%)", "12 (38%)",
    "13 (41%)", "14 (44%)", "15 (47%)", "16 (50%)", "17 (53%)", "18 (56%)",
    "19 (59%)", "20 (62%)", "21 (66%)", "22 (69%)", "23 (72%)", "24 (75%)",
    "25 (78%)", "26 (81%)", "27 (84%)", "28 (88%)", "29 (91%)", "30 (94%)",
    "31 (97%)", "32 (100%)"
)

# --- Build the target cell style once, off to the side, so only a single
#     incremental style change happens per property (this keeps the style
#     table minimal/clean instead of accumulating transient duplicates). ---
$scratch = $ws.Range("Z1")
$scratch.Value = "scratch"
$ws.Range("G3").Copy()                 # base: fontId "calibri", General number format
$scratch.PasteSpecial(-4122)           # -4122 = xlPasteFormats
$scratch.VerticalAlignment = -4107     # -4107 = xlBottom
$scratch.HorizontalAlignment = -4152   # -4152 = xlRight

# New header for column K (row 2)
$ws.Range("K2").Value = "test_xlr_n_percent"

# Fill K3:K34 with the text labels
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 3
    $ws.Cells.Item($row, 11).Value = $values[$i]
}

# Apply the prebuilt style to K3:K35 (the trailing blank table row included)
$scratch.Copy()
$ws.Range("K3:K35").PasteSpecial(-4122)

# Remove the scratch helper cell so it leaves no trace in the workbook
$scratch.Clear()

# Grow the existing table to include the new column
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A2:K34"))

# Re-assert header text (Resize can otherwise leave an autogenerated name)
$ws.Range("K2").Value = "test_xlr_n_percent"
